$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.365747
$ws.Range("H2").Value = 2.731494
$ws.Range("I2").Value = 0.01468609951202811
$ws.Range("J2").Value = 0.0101986501662014
$ws.Range("M2").Value = 181.556244
$ws.Range("N2").Value = 544.668732
$ws.Range("O2").Value = 0.393453292404907
$ws.Range("P2").Value = 0.3935455037432071
$ws.Range("Q2").Value = 247.959895574268
$ws.Range("R2").Value = 1487.759373445608
$ws.Range("S2").Value = 0.005778294205593558
$ws.Range("T2").Value = 0.004013632917158473
$ws.Range("G3").Value = 1.365747
$ws.Range("H3").Value = 2.731494
$ws.Range("I3").Value = 0.01468609951202811
$ws.Range("J3").Value = 0.0101986501662014
$ws.Range("M3").Value = 0.324361
$ws.Range("N3").Value = 0.648722
$ws.Range("O3").Value = 0.0007029276469155644
$ws.Range("P3").Value = 0.0004687282586276696
$ws.Range("Q3").Value = 0.442995062667
$ws.Range("R3").Value = 1.771980250668
$ws.Range("S3").Value = 0.00001032326537235774
$ws.Range("T3").Value = 0.000004780395532756375
$ws.Range("G4").Value = 1.365747
$ws.Range("H4").Value = 2.731494
$ws.Range("I4").Value = 0.01468609951202811
$ws.Range("J4").Value = 0.0101986501662014
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 155.929759
$ws.Range("N4").Value = 467.789277
$ws.Range("O4").Value = 0.3379177477501335
$ws.Range("P4").Value = 0.3379969435488647
$ws.Range("Q4").Value = 212.960600564973
$ws.Range("R4").Value = 1277.763603389838
$ws.Range("S4").Value = 0.004962693670338874
$ws.Range("T4").Value = 0.003447112584500194
$ws.Range("G5").Value = 1.365747
$ws.Range("H5").Value = 2.731494
$ws.Range("I5").Value = 0.01468609951202811
$ws.Range("J5").Value = 0.0101986501662014
$ws.Range("M5").Value = 123.632576
$ws.Range("N5").Value = 370.897728
$ws.Range("O5").Value = 0.2679260321980438
$ws.Range("P5").Value = 0.2679888244493004
$ws.Range("Q5").Value = 168.850819774272
$ws.Range("R5").Value = 1013.104918645632
$ws.Range("S5").Value = 0.003934788370723319
$ws.Range("T5").Value = 0.002733124269009975
$ws.Range("I6").Value = 0.2395044944124353
$ws.Range("J6").Value = 0.249483113240329
$ws.Range("M6").Value = 181.556244
$ws.Range("N6").Value = 544.668732
$ws.Range("O6").Value = 0.393453292404907
$ws.Range("P6").Value = 0.3935455037432071
$ws.Range("Q6").Value = 4043.790481974888
$ws.Range("R6").Value = 36394.11433777399
$ws.Range("S6").Value = 0.09423383187234534
$ws.Range("T6").Value = 0.09818295747558885
$ws.Range("I7").Value = 0.2395044944124353
$ws.Range("J7").Value = 0.249483113240329
$ws.Range("M7").Value = 0.324361
$ws.Range("N7").Value = 0.648722
$ws.Range("O7").Value = 0.0007029276469155644
$ws.Range("P7").Value = 0.0004687282586276696
$ws.Range("Q7").Value = 7.224471577655334
$ws.Range("R7").Value = 43.346829465932
$ws.Range("S7").Value = 0.0001683543306830351
$ws.Range("T7").Value = 0.0001169397852261491
$ws.Range("I8").Value = 0.2395044944124353
$ws.Range("J8").Value = 0.249483113240329
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 155.929759
$ws.Range("N8").Value = 467.789277
$ws.Range("O8").Value = 0.3379177477501335
$ws.Range("P8").Value = 0.3379969435488647
$ws.Range("Q8").Value = 3473.013438749252
$ws.Range("R8").Value = 31257.12094874326
$ws.Range("S8").Value = 0.08093281932788458
$ws.Range("T8").Value = 0.08432452974228648
$ws.Range("I9").Value = 0.2395044944124353
$ws.Range("J9").Value = 0.249483113240329
$ws.Range("M9").Value = 123.632576
$ws.Range("N9").Value = 370.897728
$ws.Range("O9").Value = 0.2679260321980438
$ws.Range("P9").Value = 0.2679888244493004
$ws.Range("Q9").Value = 2753.660370341419
$ws.Range("R9").Value = 24782.94333307277
$ws.Range("S9").Value = 0.06416948888152235
$ws.Range("T9").Value = 0.06685868623722745
$ws.Range("G10").Value = 36.57125933333334
$ws.Range("H10").Value = 109.713778
$ws.Range("I10").Value = 0.3932566967743808
$ws.Range("J10").Value = 0.4096411854590505
$ws.Range("M10").Value = 181.556244
$ws.Range("N10").Value = 544.668732
$ws.Range("O10").Value = 0.393453292404907
$ws.Range("P10").Value = 0.3935455037432071
$ws.Range("Q10").Value = 6639.740482909944
$ws.Range("R10").Value = 59757.66434618949
$ws.Range("S10").Value = 0.1547281421061583
$ws.Range("T10").Value = 0.1612124466854466
$ws.Range("G11").Value = 36.57125933333334
$ws.Range("H11").Value = 109.713778
$ws.Range("I11").Value = 0.3932566967743808
$ws.Range("J11").Value = 0.4096411854590505
$ws.Range("M11").Value = 0.324361
$ws.Range("N11").Value = 0.648722
$ws.Range("O11").Value = 0.0007029276469155644
$ws.Range("P11").Value = 0.0004687282586276696
$ws.Range("Q11").Value = 11.86229024861933
$ws.Range("R11").Value = 71.173741491716
$ws.Range("S11").Value = 0.0002764310044974031
$ws.Range("T11").Value = 0.000192010399522395
$ws.Range("G12").Value = 36.57125933333334
$ws.Range("H12").Value = 109.713778
$ws.Range("I12").Value = 0.3932566967743808
$ws.Range("J12").Value = 0.4096411854590505
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 155.929759
$ws.Range("N12").Value = 467.789277
$ws.Range("O12").Value = 0.3379177477501335
$ws.Range("P12").Value = 0.3379969435488647
$ws.Range("Q12").Value = 5702.547654173169
$ws.Range("R12").Value = 51322.92888755851
$ws.Range("S12").Value = 0.132888417261656
$ws.Range("T12").Value = 0.1384574686368927
$ws.Range("G13").Value = 36.57125933333334
$ws.Range("H13").Value = 109.713778
$ws.Range("I13").Value = 0.3932566967743808
$ws.Range("J13").Value = 0.4096411854590505
$ws.Range("M13").Value = 123.632576
$ws.Range("N13").Value = 370.897728
$ws.Range("O13").Value = 0.2679260321980438
$ws.Range("P13").Value = 0.2679888244493004
$ws.Range("Q13").Value = 4521.398998944043
$ws.Range("R13").Value = 40692.59099049639
$ws.Range("S13").Value = 0.1053637064020691
$ws.Range("T13").Value = 0.1097792597371888
$ws.Range("G14").Value = 9.7929715
$ws.Range("H14").Value = 19.585943
$ws.Range("I14").Value = 0.105305414522203
$ws.Range("J14").Value = 0.07312854461044438
$ws.Range("M14").Value = 181.556244
$ws.Range("N14").Value = 544.668732
$ws.Range("O14").Value = 0.393453292404907
$ws.Range("P14").Value = 0.3935455037432071
$ws.Range("Q14").Value = 1777.975123139046
$ws.Range("R14").Value = 10667.85073883428
$ws.Range("S14").Value = 0.04143276205182429
$ws.Range("T14").Value = 0.02877940992672493
$ws.Range("G15").Value = 9.7929715
$ws.Range("H15").Value = 19.585943
$ws.Range("I15").Value = 0.105305414522203
$ws.Range("J15").Value = 0.07312854461044438
$ws.Range("M15").Value = 0.324361
$ws.Range("N15").Value = 0.648722
$ws.Range("O15").Value = 0.0007029276469155644
$ws.Range("P15").Value = 0.0004687282586276696
$ws.Range("Q15").Value = 3.1764580287115
$ws.Range("R15").Value = 12.705832114846
$ws.Range("S15").Value = 0.00007402208723756026
$ws.Range("T15").Value = 0.00003427741537122945
$ws.Range("G16").Value = 9.7929715
$ws.Range("H16").Value = 19.585943
$ws.Range("I16").Value = 0.105305414522203
$ws.Range("J16").Value = 0.07312854461044438
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 155.929759
$ws.Range("N16").Value = 467.789277
$ws.Range("O16").Value = 0.3379177477501335
$ws.Range("P16").Value = 0.3379969435488647
$ws.Range("Q16").Value = 1527.015685888869
$ws.Range("R16").Value = 9162.094115333211
$ws.Range("S16").Value = 0.03558456850123704
$ws.Range("T16").Value = 0.024717224564507
$ws.Range("G17").Value = 9.7929715
$ws.Range("H17").Value = 19.585943
$ws.Range("I17").Value = 0.105305414522203
$ws.Range("J17").Value = 0.07312854461044438
$ws.Range("M17").Value = 123.632576
$ws.Range("N17").Value = 370.897728
$ws.Range("O17").Value = 0.2679260321980438
$ws.Range("P17").Value = 0.2679888244493004
$ws.Range("Q17").Value = 1210.730293239584
$ws.Range("R17").Value = 7264.381759437505
$ws.Range("S17").Value = 0.02821406188190411
$ws.Range("T17").Value = 0.01959763270384121
$ws.Range("G18").Value = 6.177412
$ws.Range("H18").Value = 18.532236
$ws.Range("I18").Value = 0.06642671545959582
$ws.Range("J18").Value = 0.06919429138833312
$ws.Range("M18").Value = 181.556244
$ws.Range("N18").Value = 544.668732
$ws.Range("O18").Value = 0.393453292404907
$ws.Range("P18").Value = 0.3935455037432071
$ws.Range("Q18").Value = 1121.547720360528
$ws.Range("R18").Value = 10093.92948324475
$ws.Range("S18").Value = 0.02613580990122191
$ws.Range("T18").Value = 0.02723110226057582
$ws.Range("G19").Value = 6.177412
$ws.Range("H19").Value = 18.532236
$ws.Range("I19").Value = 0.06642671545959582
$ws.Range("J19").Value = 0.06919429138833312
$ws.Range("M19").Value = 0.324361
$ws.Range("N19").Value = 0.648722
$ws.Range("O19").Value = 0.0007029276469155644
$ws.Range("P19").Value = 0.0004687282586276696
$ws.Range("Q19").Value = 2.003711533732
$ws.Range("R19").Value = 12.022269202392
$ws.Range("S19").Value = 0.00004669317479034344
$ws.Range("T19").Value = 0.00003243331970942894
$ws.Range("G20").Value = 6.177412
$ws.Range("H20").Value = 18.532236
$ws.Range("I20").Value = 0.06642671545959582
$ws.Range("J20").Value = 0.06919429138833312
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 155.929759
$ws.Range("N20").Value = 467.789277
$ws.Range("O20").Value = 0.3379177477501335
$ws.Range("P20").Value = 0.3379969435488647
$ws.Range("Q20").Value = 963.2423644037082
$ws.Range("R20").Value = 8669.181279633372
$ws.Range("S20").Value = 0.02244676607854559
$ws.Range("T20").Value = 0.02338745900028612
$ws.Range("G21").Value = 6.177412
$ws.Range("H21").Value = 18.532236
$ws.Range("I21").Value = 0.06642671545959582
$ws.Range("J21").Value = 0.06919429138833312
$ws.Range("M21").Value = 123.632576
$ws.Range("N21").Value = 370.897728
$ws.Range("O21").Value = 0.2679260321980438
$ws.Range("P21").Value = 0.2679888244493004
$ws.Range("Q21").Value = 763.7293585733121
$ws.Range("R21").Value = 6873.564227159809
$ws.Range("S21").Value = 0.01779744630503796
$ws.Range("T21").Value = 0.01854329680776174
$ws.Range("G22").Value = 16.815572
$ws.Range("H22").Value = 50.446716
$ws.Range("I22").Value = 0.1808205793193568
$ws.Range("J22").Value = 0.1883542151356418
$ws.Range("M22").Value = 181.556244
$ws.Range("N22").Value = 544.668732
$ws.Range("O22").Value = 0.393453292404907
$ws.Range("P22").Value = 0.3935455037432071
$ws.Range("Q22").Value = 3052.972093031568
$ws.Range("R22").Value = 27476.74883728411
$ws.Range("S22").Value = 0.07114445226776357
$ws.Range("T22").Value = 0.07412595447771257
$ws.Range("G23").Value = 16.815572
$ws.Range("H23").Value = 50.446716
$ws.Range("I23").Value = 0.1808205793193568
$ws.Range("J23").Value = 0.1883542151356418
$ws.Range("M23").Value = 0.324361
$ws.Range("N23").Value = 0.648722
$ws.Range("O23").Value = 0.0007029276469155644
$ws.Range("P23").Value = 0.0004687282586276696
$ws.Range("Q23").Value = 5.454315749492
$ws.Range("R23").Value = 32.72589449695199
$ws.Range("S23").Value = 0.0001271037843348647
$ws.Range("T23").Value = 0.00008828694326571083
$ws.Range("G24").Value = 16.815572
$ws.Range("H24").Value = 50.446716
$ws.Range("I24").Value = 0.1808205793193568
$ws.Range("J24").Value = 0.1883542151356418
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 155.929759
$ws.Range("N24").Value = 467.789277
$ws.Range("O24").Value = 0.3379177477501335
$ws.Range("P24").Value = 0.3379969435488647
$ws.Range("Q24").Value = 2622.048089407148
$ws.Range("R24").Value = 23598.43280466433
$ws.Range("S24").Value = 0.06110248291047142
$ws.Range("T24").Value = 0.06366314902039225
$ws.Range("G25").Value = 16.815572
$ws.Range("H25").Value = 50.446716
$ws.Range("I25").Value = 0.1808205793193568
$ws.Range("J25").Value = 0.1883542151356418
$ws.Range("M25").Value = 123.632576
$ws.Range("N25").Value = 370.897728
$ws.Range("O25").Value = 0.2679260321980438
$ws.Range("P25").Value = 0.2679888244493004
$ws.Range("Q25").Value = 2078.952483273472
$ws.Range("R25").Value = 18710.57234946125
$ws.Range("S25").Value = 0.04844654035678692
$ws.Range("T25").Value = 0.05047682469427128
